# Applies the updated cryptocurrency price/volume snapshot to Sheet1.
# Text-like columns (B: Coin name, C: Link) are plain strings.
# Numeric-looking columns (D: Price, E: Volume(1h)) must stay as text,
# matching the source data which stores them as strings (e.g. "320.09", "3.70%").
# Setting NumberFormat to "@" (Text) before assignment prevents Excel from
# auto-converting these strings into numeric/percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "320.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.70%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.85%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.247"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.41%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07727"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.692"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.47%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9445"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.88%"
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.58%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1241"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.08%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1829"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.16%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09213"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.34%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04343"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.77%"
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.63%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001282"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.17%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005930"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.30%"
# Row 16
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.007491"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1,897.31%"
# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.340"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.28%"
# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.339"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.35%"
# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3359"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.36%"
# Row 20
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.721"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "11.65%"
# Row 21
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1353"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.89%"
# Row 22
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2824"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.39%"
# Row 23
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04035"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.08%"
# Row 24
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001265"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.47%"
# Row 25
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004113"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.30%"
# Row 26
$ws.Range("B26").Value = "NitroEx"
$ws.Range("C26").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001271"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.11%"
# Row 27
$ws.Range("B27").Value = "Spectre.aiUtilityToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
# Row 28
$ws.Range("B28").Value = "LegolasExchange"
$ws.Range("C28").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
# Row 29
$ws.Range("B29").Value = "BitZToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
# Row 30
$ws.Range("B30").Value = "Birake"
$ws.Range("C30").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
# Row 31
$ws.Range("B31").Value = "NashExchange"
$ws.Range("C31").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
# Row 32
$ws.Range("B32").Value = "AAXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
# Row 33
$ws.Range("B33").Value = "CenX"
$ws.Range("C33").Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
# Row 34
$ws.Range("B34").Value = "BNIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02549"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.97%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05340"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.02%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007780"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.76%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1320"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.48%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007366"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.25%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001991"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.92%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008369"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "13.59%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3179"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.69%"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006674"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.24%"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.10%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2008"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "84.30%"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004203"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.99%"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.10%"
